$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# Header row (row 1): add metadata columns E:K
$ws.Range("E1").Value = "insurance"
$ws.Range("F1").Value = "normal"
$ws.Range("G1").Value = "2012-03-30"
$ws.Range("H1").Value = "黃志雄"
$ws.Range("I1").Value = 1366
$ws.Range("J1").Value = "tmpb8fa1"
$ws.Range("K1").Value = 130

# Data rows 2-17: add metadata columns E:K (category/date/legislator info + repeated index)
$indices = @(130,131,132,133,134,135,136,137,138,139,140,141,142,143,144,145)
for ($i = 0; $i -lt $indices.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = "insurance"
    $ws.Cells.Item($row, 6).Value = "normal"
    $ws.Cells.Item($row, 7).Value = "2012-03-30"
    $ws.Cells.Item($row, 8).Value = "黃志雄"
    $ws.Cells.Item($row, 9).Value = 1366
    $ws.Cells.Item($row, 10).Value = "tmpb8fa1"
    $ws.Cells.Item($row, 11).Value = $indices[$i]
}

# Remove the old "債務" (debt/investment) placeholder sheet entirely
$wsDebt = $wb.Worksheets.Item("債務")
$wsDebt.Delete()
